$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price), E (Volume 1h %) and G (Hora) in this sheet are stored
# as text (e.g. "288.59", "1.16%", "20") rather than numbers/percentages,
# so force a text number format before writing the values to preserve the
# exact string representation used by the data source.
$ws.Range("D2:D51").NumberFormat = "@"
$ws.Range("E2:E51").NumberFormat = "@"
$ws.Range("G2:G51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "288.59"
$ws.Range("E2").Value = "1.16%"
$ws.Range("G2").Value = "20"

# Row 3
$ws.Range("D3").Value = "29.15"
$ws.Range("E3").Value = "1.08%"
$ws.Range("G3").Value = "20"

# Row 4
$ws.Range("D4").Value = "5.336"
$ws.Range("E4").Value = "6.94%"
$ws.Range("G4").Value = "20"

# Row 5
$ws.Range("D5").Value = "0.07004"
$ws.Range("E5").Value = "4.59%"
$ws.Range("G5").Value = "20"

# Row 6
$ws.Range("D6").Value = "7.450"
$ws.Range("E6").Value = "1.66%"
$ws.Range("G6").Value = "20"

# Row 7
$ws.Range("D7").Value = "3.557"
$ws.Range("G7").Value = "20"

# Row 8
$ws.Range("D8").Value = "1.395"
$ws.Range("E8").Value = "2.26%"
$ws.Range("G8").Value = "20"

# Row 9
$ws.Range("D9").Value = "0.9047"
$ws.Range("E9").Value = "-3.68%"
$ws.Range("G9").Value = "20"

# Row 10
$ws.Range("D10").Value = "0.1607"
$ws.Range("E10").Value = "1.77%"
$ws.Range("G10").Value = "20"

# Row 11
$ws.Range("D11").Value = "0.07493"
$ws.Range("E11").Value = "13.90%"
$ws.Range("G11").Value = "20"

# Row 12
$ws.Range("D12").Value = "0.07723"
$ws.Range("E12").Value = "1.94%"
$ws.Range("G12").Value = "20"

# Row 13
$ws.Range("D13").Value = "0.02932"
$ws.Range("E13").Value = "-0.32%"
$ws.Range("G13").Value = "20"

# Row 14
$ws.Range("E14").Value = "0.20%"
$ws.Range("G14").Value = "20"

# Row 15
$ws.Range("D15").Value = "0.001578"
$ws.Range("E15").Value = "-1.63%"
$ws.Range("G15").Value = "20"

# Row 16
$ws.Range("E16").Value = "0.78%"
$ws.Range("G16").Value = "20"

# Row 17
$ws.Range("D17").Value = "0.006205"
$ws.Range("E17").Value = "-2.98%"
$ws.Range("G17").Value = "20"

# Row 18
$ws.Range("D18").Value = "3.477"
$ws.Range("E18").Value = "-0.47%"
$ws.Range("G18").Value = "20"

# Row 19
$ws.Range("D19").Value = "2.231"
$ws.Range("E19").Value = "-0.72%"
$ws.Range("G19").Value = "20"

# Row 20
$ws.Range("D20").Value = "0.3271"
$ws.Range("E20").Value = "2.10%"
$ws.Range("G20").Value = "20"

# Row 21
$ws.Range("D21").Value = "0.1335"
$ws.Range("E21").Value = "2.01%"
$ws.Range("G21").Value = "20"

# Row 22
$ws.Range("D22").Value = "4.007"
$ws.Range("E22").Value = "-1.39%"
$ws.Range("G22").Value = "20"

# Row 23
$ws.Range("E23").Value = "4.77%"
$ws.Range("G23").Value = "20"

# Row 24
$ws.Range("D24").Value = "0.04531"
$ws.Range("E24").Value = "0.98%"
$ws.Range("G24").Value = "20"

# Row 25
$ws.Range("D25").Value = "0.001208"
$ws.Range("E25").Value = "2.12%"
$ws.Range("G25").Value = "20"

# Row 26
$ws.Range("D26").Value = "0.004148"
$ws.Range("E26").Value = "-7.73%"
$ws.Range("G26").Value = "20"

# Row 27
$ws.Range("D27").Value = "0.0001167"
$ws.Range("E27").Value = "-6.56%"
$ws.Range("G27").Value = "20"

# Row 28
$ws.Range("D28").Value = "0.0001667"
$ws.Range("E28").Value = "3.22%"
$ws.Range("G28").Value = "20"

# Row 29
$ws.Range("G29").Value = "20"

# Row 30
$ws.Range("G30").Value = "20"

# Row 31
$ws.Range("G31").Value = "20"

# Row 32
$ws.Range("G32").Value = "20"

# Row 33
$ws.Range("G33").Value = "20"

# Row 34
$ws.Range("G34").Value = "20"

# Row 35
$ws.Range("G35").Value = "20"

# Row 36
$ws.Range("G36").Value = "20"

# Row 37
$ws.Range("G37").Value = "20"

# Row 38
$ws.Range("G38").Value = "20"

# Row 39
$ws.Range("G39").Value = "20"

# Row 40
$ws.Range("D40").Value = "0.04369"
$ws.Range("E40").Value = "3.91%"
$ws.Range("G40").Value = "20"

# Row 41
$ws.Range("D41").Value = "0.006920"
$ws.Range("E41").Value = "3.13%"
$ws.Range("G41").Value = "20"

# Row 42
$ws.Range("E42").Value = "-0.47%"
$ws.Range("G42").Value = "20"

# Row 43
$ws.Range("D43").Value = "0.002064"
$ws.Range("E43").Value = "2.29%"
$ws.Range("G43").Value = "20"

# Row 44
$ws.Range("D44").Value = "0.01165"
$ws.Range("E44").Value = "-4.41%"
$ws.Range("G44").Value = "20"

# Row 45
$ws.Range("D45").Value = "0.00005806"
$ws.Range("E45").Value = "4.44%"
$ws.Range("G45").Value = "20"

# Row 46
$ws.Range("E46").Value = "-1.85%"
$ws.Range("G46").Value = "20"

# Row 47
$ws.Range("E47").Value = "-0.37%"
$ws.Range("G47").Value = "20"

# Row 48
$ws.Range("G48").Value = "20"

# Row 49
$ws.Range("G49").Value = "20"

# Row 50
$ws.Range("G50").Value = "20"

# Row 51
$ws.Range("G51").Value = "20"

